# Update column F (dSF) values for the specified rows, per the "repull data,
# push all data, mean calculation" commit. These are the dSF figures that
# were recomputed after re-pulling the underlying data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 7
$ws.Range("F5").Value = -8
$ws.Range("F6").Value = 1
$ws.Range("F10").Value = 7
$ws.Range("F12").Value = -7
$ws.Range("F14").Value = -1
